# Applies the 4 text edits described by the diff.
# Each edit is performed by locating the target text with Find and then
# overwriting the Range.Text in place, which (for ranges that coincide
# exactly with a single run's text) preserves that run's formatting and
# does not disturb neighbouring runs.

$d = $word.ActiveDocument

# 1) "pour esbaucher, puys oings ce dedans ainsy " ->
#    "pour esbaucher, puys oings ce dedans, ainsy "
$rng = $d.Content
$found = $rng.Find.Execute("pour esbaucher, puys oings ce dedans ainsy ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) {
    $rng.Text = "pour esbaucher, puys oings ce dedans, ainsy "
}

# 2) the "</m>" immediately following "ainsy <m>ciré" becomes "</m>, "
#    (there are many "</m>" strings in the document, so first locate the
#    unique surrounding context, then narrow down to just the "</m>" run
#    so that the neighbouring "cir"/"e"-accent runs are left untouched)
$anchor = "ainsy <m>ciré</m>d"
$tagText = "</m>"
$rng2 = $d.Content
$found2 = $rng2.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false)
if ($found2) {
    $base = $rng2.Start
    $offset = $anchor.IndexOf($tagText)
    $tag = $d.Range($base + $offset, $base + $offset + $tagText.Length)
    if ($tag.Text -eq $tagText) {
        $tag.Text = "</m>, "
    }
}

# 3) "puys avecq une poincte d" -> "puys, avecq une poincte d"
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("puys avecq une poincte d", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found3) {
    $rng3.Text = "puys, avecq une poincte d"
}

# 4) " un martelet à tatiner joincts" -> " un martelet à tatiner, joincts"
$rng4 = $d.Content
$found4 = $rng4.Find.Execute(" un martelet à tatiner joincts", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found4) {
    $rng4.Text = " un martelet à tatiner, joincts"
}

Write-Output "done: $found $found2 $found3 $found4"
